# Rename the "barang_id" header (B1) to "barang_kode" on Sheet1.
# Only the header text changes; the underlying code values in B2:B6
# (SBK-004, SNK-004, MND-004, BAY-004, MNM-004) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "barang_kode"

# Column B needs to widen slightly to fit the new header text
# (matches the saved column width for the "barang_kode" header).
$ws.Columns.Item(2).ColumnWidth = 11.25

# Move the active selection to N2, matching the saved view state.
$ws.Range("N2").Select() | Out-Null
